# Update transition-probability matrix values on Sheet1 to reflect the
# refreshed "team spec time" figures (commit: "added team spec time commit pt2").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1977186311787072
$ws.Range("C2").Value = 0.5817490494296578
$ws.Range("J2").Value = 0.007604562737642586
$ws.Range("P2").Value = 0.1216730038022814
$ws.Range("S2").Value = 0.09125475285171103
$ws.Range("B3").Value = 0.01265822784810127
$ws.Range("C3").Value = 0.02531645569620253
$ws.Range("J3").Value = 0.0379746835443038
$ws.Range("P3").Value = 0.7721518987341772
$ws.Range("S3").Value = 0.1518987341772152
$ws.Range("P4").Value = 0.7608695652173914
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("B6").Value = 0.03883495145631068
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.06310679611650485
$ws.Range("J6").Value = 0.3058252427184466
$ws.Range("O6").Value = 0.01941747572815534
$ws.Range("Q6").Value = 0.1504854368932039
$ws.Range("R6").Value = 0.1019417475728155
$ws.Range("S6").Value = 0.3106796116504854
$ws.Range("B7").Value = 0.1098265895953757
$ws.Range("D7").Value = 0.005780346820809248
$ws.Range("E7").Value = 0.005780346820809248
$ws.Range("F7").Value = 0.04046242774566474
$ws.Range("J7").Value = 0.138728323699422
$ws.Range("O7").Value = 0.005780346820809248
$ws.Range("Q7").Value = 0.1849710982658959
$ws.Range("R7").Value = 0.07514450867052024
$ws.Range("S7").Value = 0.4335260115606936
$ws.Range("B8").Value = 0.09420289855072464
$ws.Range("D8").Value = 0.00966183574879227
$ws.Range("F8").Value = 0.04106280193236715
$ws.Range("J8").Value = 0.09903381642512077
$ws.Range("O8").Value = 0.02173913043478261
$ws.Range("Q8").Value = 0.2367149758454106
$ws.Range("R8").Value = 0.09903381642512077
$ws.Range("S8").Value = 0.3985507246376812
$ws.Range("B9").Value = 0.09154929577464789
$ws.Range("D9").Value = 0.01408450704225352
$ws.Range("F9").Value = 0.04225352112676056
$ws.Range("J9").Value = 0.1549295774647887
$ws.Range("O9").Value = 0.007042253521126761
$ws.Range("Q9").Value = 0.1971830985915493
$ws.Range("R9").Value = 0.1197183098591549
$ws.Range("S9").Value = 0.3732394366197183
$ws.Range("B10").Value = 0.09246088193456614
$ws.Range("D10").Value = 0.02773826458036984
$ws.Range("E10").Value = 0.001422475106685633
$ws.Range("F10").Value = 0.06258890469416785
$ws.Range("J10").Value = 0.1322901849217639
$ws.Range("O10").Value = 0.009957325746799431
$ws.Range("Q10").Value = 0.2382645803698435
$ws.Range("R10").Value = 0.07823613086770982
$ws.Range("S10").Value = 0.3570412517780939
$ws.Range("G11").Value = 0.18
$ws.Range("J11").Value = 0.09666666666666666
$ws.Range("K11").Value = 0.2266666666666667
$ws.Range("L11").Value = 0.46
$ws.Range("S11").Value = 0.03666666666666667
$ws.Range("G12").Value = 0.7132867132867133
$ws.Range("J12").Value = 0.2237762237762238
$ws.Range("K12").Value = 0.006993006993006993
$ws.Range("L12").Value = 0.01398601398601399
$ws.Range("S12").Value = 0.04195804195804196
$ws.Range("G13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.3823529411764706
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.008620689655172414
$ws.Range("H15").Value = 0.103448275862069
$ws.Range("I15").Value = 0.06465517241379311
$ws.Range("J15").Value = 0.4439655172413793
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("M15").Value = 0.004310344827586207
$ws.Range("O15").Value = 0.05603448275862069
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.01092896174863388
$ws.Range("H16").Value = 0.1311475409836066
$ws.Range("I16").Value = 0.06557377049180328
$ws.Range("J16").Value = 0.4972677595628415
$ws.Range("K16").Value = 0.08196721311475409
$ws.Range("M16").Value = 0.02185792349726776
$ws.Range("N16").Value = 0.00546448087431694
$ws.Range("O16").Value = 0.08743169398907104
$ws.Range("S16").Value = 0.09836065573770492
$ws.Range("F17").Value = 0.03461538461538462
$ws.Range("H17").Value = 0.1865384615384615
$ws.Range("I17").Value = 0.0576923076923077
$ws.Range("J17").Value = 0.4211538461538462
$ws.Range("K17").Value = 0.1019230769230769
$ws.Range("M17").Value = 0.01538461538461539
$ws.Range("O17").Value = 0.06153846153846154
$ws.Range("S17").Value = 0.1211538461538462
$ws.Range("F18").Value = 0.02985074626865672
$ws.Range("H18").Value = 0.1791044776119403
$ws.Range("I18").Value = 0.06467661691542288
$ws.Range("J18").Value = 0.4527363184079602
$ws.Range("K18").Value = 0.0945273631840796
$ws.Range("M18").Value = 0.01492537313432836
$ws.Range("O18").Value = 0.06965174129353234
$ws.Range("S18").Value = 0.0945273631840796
$ws.Range("F19").Value = 0.01761409127301842
$ws.Range("H19").Value = 0.1897518014411529
$ws.Range("I19").Value = 0.05764611689351481
$ws.Range("J19").Value = 0.4011208967173739
$ws.Range("K19").Value = 0.0976781425140112
$ws.Range("M19").Value = 0.01601281024819856
$ws.Range("O19").Value = 0.08006405124099279
$ws.Range("S19").Value = 0.1401120896717374
